# Insert a new weekly price-record row (row 294) for
# "Macroferia Regional de Talca - Acelga", pushing the existing
# rows 294-316 down to 295-317.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294 (shifts 294:316 -> 295:317,
# inheriting formatting/number-format from the row above, same as
# Excel's native "Insert Copied/Above Cells" behaviour).
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new weekly data point.
$ws.Cells.Item(294, 1).Value = 5
$ws.Cells.Item(294, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(294, 3).Value = "Maule"
$ws.Cells.Item(294, 4).Value = 44826
$ws.Cells.Item(294, 5).Value = 7
$ws.Cells.Item(294, 6).Value = 100112009
$ws.Cells.Item(294, 7).Value = "Acelga"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 500
$ws.Cells.Item(294, 11).Value = 2500
$ws.Cells.Item(294, 12).Value = 2500
$ws.Cells.Item(294, 13).Value = 2500
$ws.Cells.Item(294, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(294, 15).Value = "Región del Maule"
$ws.Cells.Item(294, 16).Value = 625
$ws.Cells.Item(294, 17).Value = 4
$ws.Cells.Item(294, 18).Value = "Hortaliza"
